$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update AO3:AO18 from 137813.87066690897 to 163644.44554783992.
# Downstream formulas in AR/AS/AT/AU recalc automatically since they
# reference AO via cell formulas.
$ws.Range("AO3:AO18").Value = 163644.44554783992
